$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 1.9
$ws.Range("O2").Value = 2

# Row 3
$ws.Range("G3").Value = 1.48
$ws.Range("I3").Value = 6.25
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 19
$ws.Range("N3").Value = 1.5
$ws.Range("O3").Value = 2.63
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.2
$ws.Range("T3").Value = 11
$ws.Range("V3").Value = 9
$ws.Range("W3").Value = 13
$ws.Range("Y3").Value = 21
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 10

# Row 4
$ws.Range("J4").Value = 1.05
$ws.Range("K4").Value = 11
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 1.8

# Row 8
$ws.Range("G8").Value = 8.5
$ws.Range("H8").Value = 4.33
$ws.Range("I8").Value = 1.38
$ws.Range("J8").Value = 1.08
$ws.Range("K8").Value = 7.5
$ws.Range("R8").Value = 2.75
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 15
$ws.Range("U8").Value = 41
$ws.Range("V8").Value = 29
$ws.Range("W8").Value = 126
$ws.Range("X8").Value = 81
$ws.Range("Y8").Value = 81
$ws.Range("AA8").Value = 9.5
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 151
$ws.Range("AE8").Value = 4.5
$ws.Range("AF8").Value = 5
$ws.Range("AH8").Value = 8

# Row 9
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 1.8
$ws.Range("X9").Value = 41
$ws.Range("AA9").Value = 6.5
$ws.Range("AB9").Value = 15
$ws.Range("AE9").Value = 7
$ws.Range("AF9").Value = 8.5
$ws.Range("AG9").Value = 8.5

# Row 11
$ws.Range("L11").Value = 1.4
$ws.Range("M11").Value = 2.5
$ws.Range("N11").Value = 2.18
$ws.Range("O11").Value = 1.53
$ws.Range("R11").Value = 2.15
$ws.Range("S11").Value = 1.55
$ws.Range("T11").Value = 5.2
$ws.Range("U11").Value = 6.5
$ws.Range("Z11").Value = 7.5
$ws.Range("AB11").Value = 22
$ws.Range("AE11").Value = 10.75
$ws.Range("AG11").Value = 18
$ws.Range("AI11").Value = 65
$ws.Range("AJ11").Value = 80

# Row 12
$ws.Range("T12").Value = 5.9
$ws.Range("U12").Value = 9.75
$ws.Range("V12").Value = 9.75
$ws.Range("X12").Value = 23
$ws.Range("Y12").Value = 45
$ws.Range("AE12").Value = 7.5
$ws.Range("AF12").Value = 15
$ws.Range("AG12").Value = 11.75
$ws.Range("AH12").Value = 40
$ws.Range("AI12").Value = 32
$ws.Range("AJ12").Value = 50

# Row 16
$ws.Range("G16").Value = 1.88
$ws.Range("H16").Value = 3.1
$ws.Range("N16").Value = 2.1
$ws.Range("O16").Value = 1.57
$ws.Range("P16").Value = 1.42
$ws.Range("Q16").Value = 2.35
$ws.Range("T16").Value = 5.3
$ws.Range("U16").Value = 7.1
$ws.Range("V16").Value = 7.1
$ws.Range("W16").Value = 13
$ws.Range("X16").Value = 13
$ws.Range("Y16").Value = 24
$ws.Range("Z16").Value = 7.5
$ws.Range("AA16").Value = 5.3
$ws.Range("AB16").Value = 13.5
$ws.Range("AC16").Value = 65
$ws.Range("AD16").Value = 500
$ws.Range("AE16").Value = 7.9
$ws.Range("AF16").Value = 16.5
$ws.Range("AG16").Value = 11.5
$ws.Range("AH16").Value = 50
$ws.Range("AI16").Value = 35
$ws.Range("AJ16").Value = 40

# Row 17
$ws.Range("G17").Value = 1.53
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 5.6
$ws.Range("N17").Value = 1.83
$ws.Range("O17").Value = 1.78
$ws.Range("P17").Value = 1.39
$ws.Range("Q17").Value = 2.45
$ws.Range("T17").Value = 5.4
$ws.Range("U17").Value = 6
$ws.Range("W17").Value = 9.25
$ws.Range("X17").Value = 10.5
$ws.Range("Y17").Value = 22
$ws.Range("Z17").Value = 9.25
$ws.Range("AA17").Value = 6.2
$ws.Range("AB17").Value = 14.5
$ws.Range("AC17").Value = 65
$ws.Range("AD17").Value = 500
$ws.Range("AE17").Value = 11.25
$ws.Range("AG17").Value = 15
$ws.Range("AH17").Value = 90
$ws.Range("AI17").Value = 50

# Row 18
$ws.Range("G18").Value = 1.38
$ws.Range("H18").Value = 4.35
$ws.Range("I18").Value = 6.5
$ws.Range("O18").Value = 2
$ws.Range("T18").Value = 6.1
$ws.Range("W18").Value = 7.5
$ws.Range("AA18").Value = 7.6
$ws.Range("AB18").Value = 16
$ws.Range("AC18").Value = 70
$ws.Range("AE18").Value = 14.5
$ws.Range("AH18").Value = 110

# Row 24
$ws.Range("J24").Value = 1.02
$ws.Range("K24").Value = 12

# Row 26
$ws.Range("G26").Value = 2
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 1.03
$ws.Range("P26").Value = 1.2
$ws.Range("Q26").Value = 4.33
$ws.Range("T26").Value = 15
$ws.Range("AA26").Value = 9
$ws.Range("AC26").Value = 26
$ws.Range("AD26").Value = 67
$ws.Range("AG26").Value = 12
$ws.Range("AJ26").Value = 21
